$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Round row 5 values down to 2 decimal places (custom accuracy) ---
$row5Values = @(15.55, 11.73, 0.4, 33.01, 27.51, 11.75, 43.95, 18.41, 8.67, 12.47, 13.44, 14.4, 4.17, 11.85, 17.17, 9.69, 0.4, 0.37, 175.6, 33.23, 11.31, 22.65, 11.72, 1.58, 21.89, 9.75, 8.68, 10.68, 14.13, 0.57, 39.72, 6.2, 13.72)

$col = 2
foreach ($v in $row5Values) {
    $ws.Cells.Item(5, $col).Value = $v
    $col = $col + 1
}

# --- 2. Delete row 6 (data trimmed to 1000 rows overall / this sheet drops its last row) ---
$ws.Rows.Item(6).Delete()

# --- 3. Narrow a subset of data columns from width 8 to width 7 ---
$narrowCols = @(2, 3, 8, 11, 12, 13, 15, 22, 24, 29, 30)
foreach ($c in $narrowCols) {
    $ws.Columns.Item($c).ColumnWidth = 7 - (5 / 6)
}
